$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pina = "Pi" + [char]0x00F1 + "a"

$ws.Rows.Item(99).Insert()

$ws.Range("A99").Value() = 4
$ws.Range("B99").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C99").Value() = "Los Lagos"
$ws.Range("D99").Value() = 44518
$ws.Range("E99").Value() = 10
$ws.Range("F99").Value() = "Fruta"
$ws.Range("G99").Value() = 100108
$ws.Range("H99").Value() = "Tropicales y subtropicales"
$ws.Range("I99").Value() = 100108005
$ws.Range("J99").Value() = $pina
$ws.Range("K99").Value() = "Caramelo"
$ws.Range("L99").Value() = "Tercera"
$ws.Range("M99").Value() = 120
$ws.Range("N99").Value() = 21000
$ws.Range("O99").Value() = 22000
$ws.Range("P99").Value() = 21500
$ws.Range("Q99").Value() = "`$/caja 16 unidades"
$ws.Range("R99").Value() = "Ecuador"
$ws.Range("S99").Value() = 1344
$ws.Range("T99").Value() = 16
